$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text storage (prevents Excel's COM
# layer from auto-coercing numeric-looking strings like "0.999" or "1.00"
# into real numbers), while preserving the cell's original style so no
# stray formatting/style gets introduced.
function Set-TextValue {
    param($addr, $val)
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = '67.454.13'
$ws.Range("E2").Value = '  -3.26%  '
$ws.Range("D3").Value = '3.705.85'
$ws.Range("E3").Value = '  -3.64%  '
Set-TextValue "D4" '0.999'
Set-TextValue "D5" '596.54'
$ws.Range("E5").Value = '  -2.16%  '
Set-TextValue "D6" '165.93'
$ws.Range("E6").Value = '  -4.96%  '
$ws.Range("D7").Value = '3.702.61'
$ws.Range("E7").Value = '  -3.57%  '
$ws.Range("E8").Value = '  -0.07%  '
Set-TextValue "D9" '0.530'
$ws.Range("E9").Value = '  +0.68%  '
Set-TextValue "D10" '0.161'
$ws.Range("E10").Value = '  -2.90%  '
$ws.Range("E11").Value = '  -4.36%  '
Set-TextValue "D12" '0.463'
$ws.Range("E12").Value = '  -3.55%  '
Set-TextValue "D13" '37.70'
$ws.Range("E13").Value = '  -5.41%  '
Set-TextValue "D14" '0.0000242'
$ws.Range("E14").Value = '  -4.68%  '
$ws.Range("D15").Value = '4.318.94'
$ws.Range("E15").Value = '  -3.71%  '
$ws.Range("D16").Value = '3.701.65'
$ws.Range("E16").Value = '  -3.48%  '
$ws.Range("D17").Value = '67.484.13'
$ws.Range("E17").Value = '  -3.30%  '
Set-TextValue "D18" '17.59'
$ws.Range("E18").Value = '  +5.93%  '
Set-TextValue "D19" '7.18'
$ws.Range("E19").Value = '  -3.74%  '
$ws.Range("E20").Value = '  -3.11%  '
Set-TextValue "D21" '492.32'
$ws.Range("E22").Value = '  -3.87%  '
$ws.Range("E23").Value = '  -1.78%  '
Set-TextValue "D24" '85.87'
$ws.Range("E24").Value = '  +0.05%  '
Set-TextValue "D25" '2.31'
$ws.Range("E25").Value = '  -5.93%  '
$ws.Range("E26").Value = '  -2.25%  '
Set-TextValue "D27" '12.18'
$ws.Range("E27").Value = '  -3.22%  '
Set-TextValue "D28" '10.12'
$ws.Range("E28").Value = '  -2.94%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  -1.66%  '
Set-TextValue "D31" '2.36'
$ws.Range("E31").Value = '  -6.27%  '
Set-TextValue "D32" '31.56'
$ws.Range("E32").Value = '  -2.53%  '
Set-TextValue "D33" '7.64'
$ws.Range("E33").Value = '  -3.62%  '
$ws.Range("D34").Value = '3.840.08'
$ws.Range("E34").Value = '  -3.69%  '
Set-TextValue "D35" '0.108'
$ws.Range("E35").Value = '  -4.29%  '
$ws.Range("D36").Value = '3.641.31'
$ws.Range("E36").Value = '  -3.75%  '
$ws.Range("E37").Value = '  +0.06%  '
Set-TextValue "D38" '0.998'
$ws.Range("E38").Value = '  -4.16%  '
$ws.Range("E39").Value = '  -5.39%  '
$ws.Range("E40").Value = '  -6.39%  '
Set-TextValue "D41" '0.323'
$ws.Range("E41").Value = '  -3.63%  '
Set-TextValue "D42" '433.79'
$ws.Range("E42").Value = '  -10.08%  '
$ws.Range("E43").Value = '  -2.23%  '
Set-TextValue "D44" '1.94'
$ws.Range("E44").Value = '  -5.41%  '
Set-TextValue "D45" '2.80'
$ws.Range("E45").Value = '  -6.38%  '
Set-TextValue "D46" '8.40'
$ws.Range("E46").Value = '  -1.18%  '
$ws.Range("B47").Value = 'Arweave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue "D47" '40.80'
$ws.Range("E47").Value = '  -5.34%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue "D48" '1.00'
$ws.Range("E48").Value = '  -0.01%  '
Set-TextValue "D49" '143.36'
$ws.Range("E49").Value = '  +2.36%  '
$ws.Range("D50").Value = '2.753.46'
$ws.Range("E50").Value = '  -5.54%  '
Set-TextValue "D51" '0.0348'
$ws.Range("E51").Value = '  -3.46%  '
